$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-87 down to 66-88.
$ws.Rows("65:65").Insert()

# Populate the newly inserted row 65 with the new weekly price record.
$ws.Range("A65").Value = 11
$ws.Range("B65").Value = "Vega Monumental Concepción"
$ws.Range("C65").Value = "Bíobío"
$ws.Range("D65").Value = 44523
$ws.Range("E65").Value = 8
$ws.Range("F65").Value = 100112043
$ws.Range("G65").Value = "Pepino ensalada"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 100
$ws.Range("K65").Value = 6500
$ws.Range("L65").Value = 7000
$ws.Range("M65").Value = 6750
$ws.Range("N65").Value = "`$/caja 60 unidades"
$ws.Range("O65").Value = "Región de Arica y Parinacota"
$ws.Range("P65").Value = 112
$ws.Range("Q65").Value = 60
$ws.Range("R65").Value = "Hortaliza"
